$wb = $excel.ActiveWorkbook

# --- "config" sheet: add a commodity column, drop the year column/rows ---
$ws = $wb.Worksheets.Item("config")
$ws.Activate()

# The sheet used to be node/sector/level/commodity/year; it becomes
# node/sector/commodity/level (commodity moves into column C, level shifts
# to D, and the year column/extra rows go away).
$ws.Range("C1").Value = "commodity"
$ws.Range("C2").Value = "light"
$ws.Range("D1").Value = "level"
$ws.Range("D2").Value = "useful"

# Remove the now-unused "year" column and the rows that only existed for it.
$ws.Range("E1:E4").ClearContents()
$ws.Range("A3:D4").ClearContents()

# Widen the new commodity column to fit its contents.
$ws.Columns.Item(3).ColumnWidth = 10.166666666666666

$ws.Range("A3").Select()

# --- "MERtoPPP" sheet: selection/view state only ---
$ws2 = $wb.Worksheets.Item("MERtoPPP")
$ws2.Activate()
$ws2.Range("K2").Select()

$ws.Activate()
